$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.7994676808118871
$ws.Range("J2").Value = 0.7994676808118871
$ws.Range("M2").Value = 10.70375566666667
$ws.Range("N2").Value = 32.111267
$ws.Range("O2").Value = 0.3178747217938744
$ws.Range("P2").Value = 0.3178747217938744
$ws.Range("Q2").Value = 2.131613693912556
$ws.Range("R2").Value = 19.184523245213
$ws.Range("S2").Value = 0.2541305666212726
$ws.Range("T2").Value = 0.2541305666212726

# Row 3
$ws.Range("I3").Value = 0.7994676808118871
$ws.Range("J3").Value = 0.7994676808118871
$ws.Range("O3").Value = 0.4731463873433873
$ws.Range("P3").Value = 0.4731463873433873
$ws.Range("S3").Value = 0.3782652449739407
$ws.Range("T3").Value = 0.3782652449739407

# Row 4
$ws.Range("I4").Value = 0.7994676808118871
$ws.Range("J4").Value = 0.7994676808118871
$ws.Range("M4").Value = 7.036919999999999
$ws.Range("N4").Value = 21.11076
$ws.Range("O4").Value = 0.2089788908627383
$ws.Range("P4").Value = 0.2089788908627384
$ws.Range("Q4").Value = 1.40137681596
$ws.Range("R4").Value = 12.61239134364
$ws.Range("S4").Value = 0.1670718692166739
$ws.Range("T4").Value = 0.1670718692166739

# Row 5
$ws.Range("G5").Value = 0.04995233333333333
$ws.Range("H5").Value = 0.149857
$ws.Range("I5").Value = 0.2005323191881128
$ws.Range("J5").Value = 0.2005323191881128
$ws.Range("M5").Value = 10.70375566666667
$ws.Range("N5").Value = 32.111267
$ws.Range("O5").Value = 0.3178747217938744
$ws.Range("P5").Value = 0.3178747217938744
$ws.Range("Q5").Value = 0.5346775709798888
$ws.Range("R5").Value = 4.812098138819
$ws.Range("S5").Value = 0.06374415517260179
$ws.Range("T5").Value = 0.06374415517260179

# Row 6
$ws.Range("G6").Value = 0.04995233333333333
$ws.Range("H6").Value = 0.149857
$ws.Range("I6").Value = 0.2005323191881128
$ws.Range("J6").Value = 0.2005323191881128
$ws.Range("O6").Value = 0.4731463873433873
$ws.Range("P6").Value = 0.4731463873433873
$ws.Range("Q6").Value = 0.7958505151809998
$ws.Range("R6").Value = 7.162654636628999
$ws.Range("S6").Value = 0.09488114236944661
$ws.Range("T6").Value = 0.09488114236944663

# Row 7
$ws.Range("G7").Value = 0.04995233333333333
$ws.Range("H7").Value = 0.149857
$ws.Range("I7").Value = 0.2005323191881128
$ws.Range("J7").Value = 0.2005323191881128
$ws.Range("M7").Value = 7.036919999999999
$ws.Range("N7").Value = 21.11076
$ws.Range("O7").Value = 0.2089788908627383
$ws.Range("P7").Value = 0.2089788908627384
$ws.Range("Q7").Value = 0.3515105734799999
$ws.Range("R7").Value = 3.16359516132
$ws.Range("S7").Value = 0.04190702164606444
$ws.Range("T7").Value = 0.04190702164606445
